$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1120.7693
$ws.Range("J17").Value = 925.82355
$ws.Range("L17").Value = 2777.47065
$ws.Range("N17").Value = -3113.47065
$ws.Range("H18").Value = 13417.533
$ws.Range("I18").Value = 3470
$ws.Range("J18").Value = 17034.818
$ws.Range("K18").Value = 3470
$ws.Range("L18").Value = 17034.818
$ws.Range("M18").Value = -3186
$ws.Range("N18").Value = -17602.818
$ws.Range("H86").Value = 1648.75
$ws.Range("I86").Value = 1600
$ws.Range("J86").Value = 1795
$ws.Range("K86").Value = 1600
$ws.Range("L86").Value = 1795
$ws.Range("M86").Value = -477
$ws.Range("N86").Value = -4041
$ws.Range("H88").Value = 3826.0667
$ws.Range("I88").Value = 4497.3335
$ws.Range("J88").Value = 3658.25
$ws.Range("K88").Value = 4497.3335
$ws.Range("L88").Value = 3658.25
$ws.Range("M88").Value = -4091.3335
$ws.Range("N88").Value = -4470.25
$ws.Range("H89").Value = 1648.75
$ws.Range("I89").Value = 1600
$ws.Range("J89").Value = 1795
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 8975
$ws.Range("M89").Value = -2384
$ws.Range("N89").Value = -20207
$ws.Range("H91").Value = 3826.0667
$ws.Range("I91").Value = 4497.3335
$ws.Range("J91").Value = 3658.25
$ws.Range("K91").Value = 4497.3335
$ws.Range("L91").Value = 3658.25
$ws.Range("M91").Value = -3093.3335
$ws.Range("N91").Value = -6466.25
$ws.Range("H112").Value = 4676.077
$ws.Range("J112").Value = 5032.4165
$ws.Range("L112").Value = 15097.2495
$ws.Range("N112").Value = -17313.2495
$ws.Range("H137").Value = 1275.6428
$ws.Range("I137").Value = 1241.7273
$ws.Range("J137").Value = 1400
$ws.Range("K137").Value = 3725.1819
$ws.Range("L137").Value = 4200
$ws.Range("M137").Value = -1175.1819
$ws.Range("N137").Value = -9300
$ws.Range("H138").Value = 2811.0588
$ws.Range("I138").Value = 2728.2307
$ws.Range("J138").Value = 2897.2
$ws.Range("K138").Value = 8184.6921
$ws.Range("L138").Value = 8691.599999999999
$ws.Range("M138").Value = -3044.6921
$ws.Range("N138").Value = -18971.6
$ws.Range("H141").Value = 1870919.6
$ws.Range("I141").Value = 3114412
$ws.Range("K141").Value = 9343236
$ws.Range("M141").Value = -9338056

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3795.6863
$ws.Range("I32").Value = 2686.7273
$ws.Range("K32").Value = 2686.7273
$ws.Range("M32").Value = -2399.7273
$ws.Range("H45").Value = 1708.2142
$ws.Range("I45").Value = 1487.3334
$ws.Range("J45").Value = 1873.875
$ws.Range("K45").Value = 1487.3334
$ws.Range("L45").Value = 1873.875
$ws.Range("M45").Value = -1110.3334
$ws.Range("N45").Value = -2627.875
$ws.Range("I88").Value = 2249.5
$ws.Range("J88").Value = 4249.5
$ws.Range("K88").Value = 2249.5
$ws.Range("L88").Value = 4249.5
$ws.Range("M88").Value = -1843.5
$ws.Range("N88").Value = -5061.5
$ws.Range("I91").Value = 2249.5
$ws.Range("J91").Value = 4249.5
$ws.Range("K91").Value = 2249.5
$ws.Range("L91").Value = 4249.5
$ws.Range("M91").Value = -845.5
$ws.Range("N91").Value = -7057.5
$ws.Range("H102").Value = 2410.4666
$ws.Range("I102").Value = 2337.3333
$ws.Range("J102").Value = 2703
$ws.Range("K102").Value = 2337.3333
$ws.Range("L102").Value = 2703
$ws.Range("M102").Value = -715.3332999999998
$ws.Range("N102").Value = -5947

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2603.8572
$ws.Range("I107").Value = 2834.0908
$ws.Range("J107").Value = 1759.6666
$ws.Range("K107").Value = 2834.0908
$ws.Range("L107").Value = 1759.6666
$ws.Range("M107").Value = -914.0907999999999
$ws.Range("N107").Value = -5599.6666
$ws.Range("H134").Value = 12446.091
$ws.Range("I134").Value = 12727.556
$ws.Range("J134").Value = 11179.5
$ws.Range("K134").Value = 38182.66800000001
$ws.Range("L134").Value = 33538.5
$ws.Range("M134").Value = -35647.66800000001
$ws.Range("N134").Value = -38608.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 797.8182
$ws.Range("I16").Value = 753.1111
$ws.Range("K16").Value = 753.1111
$ws.Range("M16").Value = -466.1111
$ws.Range("H31").Value = 3033.606
$ws.Range("I31").Value = 2076.25
$ws.Range("J31").Value = 4506.4614
$ws.Range("K31").Value = 2076.25
$ws.Range("L31").Value = 4506.4614
$ws.Range("M31").Value = -1781.25
$ws.Range("N31").Value = -5096.4614
$ws.Range("H34").Value = 3033.606
$ws.Range("I34").Value = 2076.25
$ws.Range("J34").Value = 4506.4614
$ws.Range("K34").Value = 2076.25
$ws.Range("L34").Value = 4506.4614
$ws.Range("M34").Value = -1874.25
$ws.Range("N34").Value = -4910.4614
$ws.Range("H94").Value = 1331
$ws.Range("J94").Value = 1331
$ws.Range("L94").Value = 1331
$ws.Range("N94").Value = -2233
$ws.Range("H113").Value = 797.8182
$ws.Range("I113").Value = 753.1111
$ws.Range("K113").Value = 753.1111
$ws.Range("M113").Value = 1416.8889
$ws.Range("H132").Value = 1768.5714
$ws.Range("I132").Value = 1017.34784
$ws.Range("J132").Value = 3208.4167
$ws.Range("K132").Value = 3052.04352
$ws.Range("L132").Value = 9625.250100000001
$ws.Range("M132").Value = -522.0435200000002
$ws.Range("N132").Value = -14685.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1309.75
$ws.Range("J34").Value = 1663
$ws.Range("L34").Value = 4989
$ws.Range("N34").Value = -5157
$ws.Range("H107").Value = 889.73334
$ws.Range("J107").Value = 935.4286
$ws.Range("L107").Value = 2806.2858
$ws.Range("N107").Value = -6646.2858
$ws.Range("H131").Value = 12743.269
$ws.Range("I131").Value = 507.8
$ws.Range("J131").Value = 14889.842
$ws.Range("K131").Value = 1523.4
$ws.Range("L131").Value = 44669.526
$ws.Range("M131").Value = 3516.6
$ws.Range("N131").Value = -54749.526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 120.21429
$ws.Range("I2").Value = 195.66667
$ws.Range("J2").Value = 63.625
$ws.Range("K2").Value = 195.66667
$ws.Range("L2").Value = 63.625
$ws.Range("M2").Value = -82.66667000000001
$ws.Range("N2").Value = -289.625
$ws.Range("H70").Value = 15785.571
$ws.Range("I70").Value = 20099.8
$ws.Range("K70").Value = 20099.8
$ws.Range("M70").Value = -19829.8
$ws.Range("H73").Value = 15785.571
$ws.Range("I73").Value = 20099.8
$ws.Range("K73").Value = 20099.8
$ws.Range("M73").Value = -19163.8
$ws.Range("H102").Value = 2920.3076
$ws.Range("I102").Value = 2950
$ws.Range("K102").Value = 2950
$ws.Range("M102").Value = -1328
$ws.Range("H107").Value = 1067.6666
$ws.Range("I107").Value = 100
$ws.Range("K107").Value = 100
$ws.Range("M107").Value = 1820
$ws.Range("H113").Value = 1366.6666
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070
$ws.Range("H122").Value = 1678.6
$ws.Range("I122").Value = 1313.4445
$ws.Range("J122").Value = 2226.3333
$ws.Range("K122").Value = 3940.3335
$ws.Range("L122").Value = 6678.999899999999
$ws.Range("M122").Value = -1490.3335
$ws.Range("N122").Value = -11578.9999
$ws.Range("H126").Value = 1826934.1
$ws.Range("I126").Value = 2418183.8
$ws.Range("K126").Value = 7254551.399999999
$ws.Range("M126").Value = -7252081.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8540.637000000001
$ws.Range("I40").Value = 3158.3333
$ws.Range("K40").Value = 3158.3333
$ws.Range("M40").Value = -3022.3333
$ws.Range("H122").Value = 7190.4
$ws.Range("I122").Value = 5238
$ws.Range("K122").Value = 15714
$ws.Range("M122").Value = -13264

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2883.3914
$ws.Range("I132").Value = 2555
$ws.Range("J132").Value = 3241.6365
$ws.Range("K132").Value = 7665
$ws.Range("L132").Value = 9724.9095
$ws.Range("M132").Value = -5135
$ws.Range("N132").Value = -14784.9095
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120
